# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new F value. F22 differs slightly between sheets in the source data
# (177 on 展览, 178 on 全部类型) but both converge to 181 after the update.
$commonUpdates = @{
    2  = 8869
    3  = 8302
    4  = 144
    12 = 754
    13 = 207
    14 = 5375
    15 = 5
    18 = 16
    21 = 158
    22 = 181
    23 = 12
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $commonUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $commonUpdates[$row]
    }
}
